# Scheduled runner update: refresh computed market-price / profit figures
# across the per-job Leve profit sheets (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H-N) with newly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 639.5
$ws.Cells.Item(2, 10).Value = 950
$ws.Cells.Item(2, 12).Value = 950
$ws.Cells.Item(2, 14).Value = -1176
$ws.Cells.Item(29, 8).Value = 337.33334
$ws.Cells.Item(29, 9).Value = 337.33334
$ws.Cells.Item(29, 11).Value = 1012.00002
$ws.Cells.Item(29, 13).Value = -731.0000200000001
$ws.Cells.Item(53, 8).Value = 859.5
$ws.Cells.Item(53, 9).Value = 1083.7142
$ws.Cells.Item(53, 10).Value = 336.33334
$ws.Cells.Item(53, 11).Value = 1083.7142
$ws.Cells.Item(53, 12).Value = 336.33334
$ws.Cells.Item(53, 13).Value = -446.7141999999999
$ws.Cells.Item(53, 14).Value = -1610.33334
$ws.Cells.Item(57, 8).Value = 48000
$ws.Cells.Item(57, 10).Value = 48000
$ws.Cells.Item(57, 12).Value = 144000
$ws.Cells.Item(57, 14).Value = -144998
$ws.Cells.Item(61, 8).Value = 999
$ws.Cells.Item(61, 9).Value = 999
$ws.Cells.Item(61, 11).Value = 2997
$ws.Cells.Item(61, 13).Value = -2825
$ws.Cells.Item(98, 8).Value = 784.875
$ws.Cells.Item(98, 9).Value = 747.7143
$ws.Cells.Item(98, 11).Value = 747.7143
$ws.Cells.Item(98, 13).Value = 750.2857
$ws.Cells.Item(107, 8).Value = 456.5
$ws.Cells.Item(107, 9).Value = 456.5
$ws.Cells.Item(107, 11).Value = 456.5
$ws.Cells.Item(107, 13).Value = 1463.5
$ws.Cells.Item(110, 8).Value = 35000
$ws.Cells.Item(110, 10).Value = 35000
$ws.Cells.Item(110, 12).Value = 35000
$ws.Cells.Item(110, 14).Value = -43180
$ws.Cells.Item(113, 8).Value = 3000
$ws.Cells.Item(113, 9).Value = 2602.5
$ws.Cells.Item(113, 11).Value = 2602.5
$ws.Cells.Item(113, 13).Value = 651.5
$ws.Cells.Item(122, 8).Value = 784.875
$ws.Cells.Item(122, 9).Value = 747.7143
$ws.Cells.Item(122, 11).Value = 2243.1429
$ws.Cells.Item(122, 13).Value = 206.8571000000002
$ws.Cells.Item(125, 8).Value = 3176.2856
$ws.Cells.Item(125, 9).Value = 1199.5
$ws.Cells.Item(125, 10).Value = 3967
$ws.Cells.Item(125, 11).Value = 10795.5
$ws.Cells.Item(125, 12).Value = 35703
$ws.Cells.Item(125, 13).Value = -8335.5
$ws.Cells.Item(125, 14).Value = -40623
$ws.Cells.Item(137, 8).Value = 1539.6364
$ws.Cells.Item(137, 9).Value = 1373.3125
$ws.Cells.Item(137, 11).Value = 4119.9375
$ws.Cells.Item(137, 13).Value = -1569.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6980.1113
$ws.Cells.Item(61, 9).Value = 7227.625
$ws.Cells.Item(61, 11).Value = 7227.625
$ws.Cells.Item(61, 13).Value = -7015.625
$ws.Cells.Item(63, 8).Value = 2000.7142
$ws.Cells.Item(63, 9).Value = 1500.8334
$ws.Cells.Item(63, 10).Value = 5000
$ws.Cells.Item(63, 11).Value = 1500.8334
$ws.Cells.Item(63, 12).Value = 5000
$ws.Cells.Item(63, 13).Value = -814.8334
$ws.Cells.Item(63, 14).Value = -6372
$ws.Cells.Item(66, 8).Value = 2000.7142
$ws.Cells.Item(66, 9).Value = 1500.8334
$ws.Cells.Item(66, 10).Value = 5000
$ws.Cells.Item(66, 11).Value = 7504.166999999999
$ws.Cells.Item(66, 12).Value = 25000
$ws.Cells.Item(66, 13).Value = -4072.166999999999
$ws.Cells.Item(66, 14).Value = -31864
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 13).ClearContents()
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 2232.7856
$ws.Cells.Item(102, 9).Value = 926
$ws.Cells.Item(102, 11).Value = 926
$ws.Cells.Item(102, 13).Value = 696
$ws.Cells.Item(136, 8).Value = 6980.1113
$ws.Cells.Item(136, 9).Value = 7227.625
$ws.Cells.Item(136, 11).Value = 21682.875
$ws.Cells.Item(136, 13).Value = -19132.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(74, 8).Value = 55388
$ws.Cells.Item(74, 10).Value = 55388
$ws.Cells.Item(74, 12).Value = 55388
$ws.Cells.Item(74, 14).Value = -57260
$ws.Cells.Item(77, 8).Value = 55388
$ws.Cells.Item(77, 10).Value = 55388
$ws.Cells.Item(77, 12).Value = 166164
$ws.Cells.Item(77, 14).Value = -175524
$ws.Cells.Item(92, 8).Value = 54999
$ws.Cells.Item(92, 10).Value = 54999
$ws.Cells.Item(92, 12).Value = 54999
$ws.Cells.Item(92, 14).Value = -59991
$ws.Cells.Item(105, 8).Value = 4494.75
$ws.Cells.Item(105, 9).Value = 4326.3335
$ws.Cells.Item(105, 11).Value = 4326.3335
$ws.Cells.Item(105, 13).Value = -2579.3335
$ws.Cells.Item(134, 8).Value = 3524.1428
$ws.Cells.Item(134, 9).Value = 1441.3334
$ws.Cells.Item(134, 11).Value = 4324.0002
$ws.Cells.Item(134, 13).Value = -1789.0002
$ws.Cells.Item(139, 8).Value = 74998.5
$ws.Cells.Item(139, 10).Value = 74998.5
$ws.Cells.Item(139, 12).Value = 74998.5
$ws.Cells.Item(139, 14).Value = -85278.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3089.889
$ws.Cells.Item(58, 9).Value = 3089.889
$ws.Cells.Item(58, 11).Value = 3089.889
$ws.Cells.Item(58, 13).Value = -2886.889
$ws.Cells.Item(92, 8).Value = 23050
$ws.Cells.Item(92, 10).Value = 23050
$ws.Cells.Item(92, 12).Value = 23050
$ws.Cells.Item(92, 14).Value = -28042
$ws.Cells.Item(99, 8).Value = 5505.75
$ws.Cells.Item(99, 9).Value = 5232.25
$ws.Cells.Item(99, 10).Value = 6326.25
$ws.Cells.Item(99, 11).Value = 5232.25
$ws.Cells.Item(99, 12).Value = 6326.25
$ws.Cells.Item(99, 13).Value = -3734.25
$ws.Cells.Item(99, 14).Value = -9322.25
$ws.Cells.Item(105, 8).Value = 1454.9333
$ws.Cells.Item(105, 9).Value = 961.3
$ws.Cells.Item(105, 11).Value = 961.3
$ws.Cells.Item(105, 13).Value = 785.7
$ws.Cells.Item(126, 8).Value = 5505.75
$ws.Cells.Item(126, 9).Value = 5232.25
$ws.Cells.Item(126, 10).Value = 6326.25
$ws.Cells.Item(126, 11).Value = 15696.75
$ws.Cells.Item(126, 12).Value = 18978.75
$ws.Cells.Item(126, 13).Value = -13226.75
$ws.Cells.Item(126, 14).Value = -23918.75
$ws.Cells.Item(136, 8).Value = 3089.889
$ws.Cells.Item(136, 9).Value = 3089.889
$ws.Cells.Item(136, 11).Value = 9269.667000000001
$ws.Cells.Item(136, 13).Value = -6719.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 9).Value = 0
$ws.Cells.Item(54, 11).Value = 0
$ws.Cells.Item(54, 13).ClearContents()
$ws.Cells.Item(68, 8).Value = 1999
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(69, 8).Value = 899
$ws.Cells.Item(69, 9).Value = 899
$ws.Cells.Item(69, 11).Value = 2697
$ws.Cells.Item(69, 13).Value = -1886
$ws.Cells.Item(71, 8).Value = 1999
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).ClearContents()
$ws.Cells.Item(72, 8).Value = 899
$ws.Cells.Item(72, 9).Value = 899
$ws.Cells.Item(72, 11).Value = 8091
$ws.Cells.Item(72, 13).Value = -4035
$ws.Cells.Item(122, 8).Value = 1976.5
$ws.Cells.Item(122, 9).Value = 1976
$ws.Cells.Item(122, 11).Value = 17784
$ws.Cells.Item(122, 13).Value = -15334

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 27269.572
$ws.Cells.Item(136, 10).Value = 27269.572
$ws.Cells.Item(136, 12).Value = 81808.716
$ws.Cells.Item(136, 14).Value = -86908.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 647.9524
$ws.Cells.Item(55, 9).Value = 165
$ws.Cells.Item(55, 11).Value = 165
$ws.Cells.Item(55, 13).Value = 8
$ws.Cells.Item(74, 8).Value = 47499.5
$ws.Cells.Item(74, 10).Value = 47499.5
$ws.Cells.Item(74, 12).Value = 47499.5
$ws.Cells.Item(74, 14).Value = -49495.5
$ws.Cells.Item(77, 8).Value = 47499.5
$ws.Cells.Item(77, 10).Value = 47499.5
$ws.Cells.Item(77, 12).Value = 142498.5
$ws.Cells.Item(77, 14).Value = -152482.5
$ws.Cells.Item(82, 8).Value = 4408.273
$ws.Cells.Item(82, 9).Value = 947.5
$ws.Cells.Item(82, 10).Value = 5177.3335
$ws.Cells.Item(82, 11).Value = 947.5
$ws.Cells.Item(82, 12).Value = 5177.3335
$ws.Cells.Item(82, 13).Value = -586.5
$ws.Cells.Item(82, 14).Value = -5899.3335
$ws.Cells.Item(85, 8).Value = 4408.273
$ws.Cells.Item(85, 9).Value = 947.5
$ws.Cells.Item(85, 10).Value = 5177.3335
$ws.Cells.Item(85, 11).Value = 947.5
$ws.Cells.Item(85, 12).Value = 5177.3335
$ws.Cells.Item(85, 13).Value = 300.5
$ws.Cells.Item(85, 14).Value = -7673.3335
$ws.Cells.Item(101, 8).Value = 27559.334
$ws.Cells.Item(101, 10).Value = 27559.334
$ws.Cells.Item(101, 12).Value = 27559.334
$ws.Cells.Item(101, 14).Value = -34049.334
$ws.Cells.Item(132, 8).Value = 4666.3335
$ws.Cells.Item(132, 9).Value = 4666.3335
$ws.Cells.Item(132, 11).Value = 13999.0005
$ws.Cells.Item(132, 13).Value = -11469.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 5810730
$ws.Cells.Item(100, 10).Value = 3254.7778
$ws.Cells.Item(100, 12).Value = 6509.5556
$ws.Cells.Item(100, 14).Value = -7591.5556
$ws.Cells.Item(112, 8).Value = 34051.332
$ws.Cells.Item(112, 10).Value = 34051.332
$ws.Cells.Item(112, 12).Value = 34051.332
$ws.Cells.Item(112, 14).Value = -37005.332
$ws.Cells.Item(125, 8).Value = 60000
$ws.Cells.Item(125, 10).Value = 60000
$ws.Cells.Item(125, 12).Value = 60000
$ws.Cells.Item(125, 14).Value = -69840
$ws.Cells.Item(133, 8).Value = 49999.5
$ws.Cells.Item(133, 10).Value = 49999.5
$ws.Cells.Item(133, 12).Value = 49999.5
$ws.Cells.Item(133, 14).Value = -60119.5
